# Penalty Reward System (unfinished): append one more weekly bucket
# (the upcoming week, currently at zero) to the two weekly-aggregated
# sheets so the forecast/PO comparison keeps pace with the new date.
#
# "Weekly Sales"        : A1:B24 -> A1:B25  (new row 25)
# "Merged (Optional)"   : A1:C32 -> A1:C33  (new row 33)
#
# Both new rows reuse the existing weekly cadence (prior date + 7 days)
# and the same date/time number format already applied to column A.

$wb = $excel.ActiveWorkbook

$newDate = 45662.99999999999

$wsWeeklySales = $wb.Worksheets.Item("Weekly Sales")
$wsWeeklySales.Cells.Item(25, 1).Value = $newDate
$wsWeeklySales.Cells.Item(25, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsWeeklySales.Cells.Item(25, 2).Value = 0

$wsMerged = $wb.Worksheets.Item("Merged (Optional)")
$wsMerged.Cells.Item(33, 1).Value = $newDate
$wsMerged.Cells.Item(33, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$wsMerged.Cells.Item(33, 2).Value = 0
$wsMerged.Cells.Item(33, 3).Value = 0
